$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 73.38544233333333
$ws.Range("H2").Value = 220.156327
$ws.Range("I2").Value = 0.1214979676060253
$ws.Range("J2").Value = 0.1214979676060253
$ws.Range("M2").Value = 15.47987166666667
$ws.Range("N2").Value = 46.439615
$ws.Range("O2").Value = 0.7960757698994193
$ws.Range("P2").Value = 0.7960757698994194
$ws.Range("Q2").Value = 1135.997229521567
$ws.Range("R2").Value = 10223.97506569411
$ws.Range("S2").Value = 0.09672158810318129
$ws.Range("T2").Value = 0.09672158810318129
$ws.Range("G3").Value = 73.38544233333333
$ws.Range("H3").Value = 220.156327
$ws.Range("I3").Value = 0.1214979676060253
$ws.Range("J3").Value = 0.1214979676060253
$ws.Range("O3").Value = 0.05100527512565552
$ws.Range("P3").Value = 0.05100527512565553
$ws.Range("Q3").Value = 72.784342175181
$ws.Range("R3").Value = 655.059079576629
$ws.Range("S3").Value = 0.006197037264953303
$ws.Range("T3").Value = 0.006197037264953303
$ws.Range("G4").Value = 73.38544233333333
$ws.Range("H4").Value = 220.156327
$ws.Range("I4").Value = 0.1214979676060253
$ws.Range("J4").Value = 0.1214979676060253
$ws.Range("M4").Value = 0.171678
$ws.Range("N4").Value = 0.515034
$ws.Range("O4").Value = 0.008828800326496623
$ws.Range("P4").Value = 0.008828800326496624
$ws.Range("Q4").Value = 12.598665968902
$ws.Range("R4").Value = 113.387993720118
$ws.Range("S4").Value = 0.001072681296068752
$ws.Range("T4").Value = 0.001072681296068752
$ws.Range("G5").Value = 73.38544233333333
$ws.Range("H5").Value = 220.156327
$ws.Range("I5").Value = 0.1214979676060253
$ws.Range("J5").Value = 0.1214979676060253
$ws.Range("M5").Value = 2.688466333333334
$ws.Range("N5").Value = 8.065399000000001
$ws.Range("O5").Value = 0.1382584398787761
$ws.Range("P5").Value = 0.1382584398787761
$ws.Range("Q5").Value = 197.2942910699415
$ws.Range("R5").Value = 1775.648619629473
$ws.Range("S5").Value = 0.01679811944965113
$ws.Range("T5").Value = 0.01679811944965113
$ws.Range("G6").Value = 73.38544233333333
$ws.Range("H6").Value = 220.156327
$ws.Range("I6").Value = 0.1214979676060253
$ws.Range("J6").Value = 0.1214979676060253
$ws.Range("M6").Value = 0.113399
$ws.Range("N6").Value = 0.340197
$ws.Range("O6").Value = 0.005831714769652435
$ws.Range("P6").Value = 0.005831714769652436
$ws.Range("Q6").Value = 8.321835775157668
$ws.Range("R6").Value = 74.89652197641901
$ws.Range("S6").Value = 0.0007085414921708108
$ws.Range("T6").Value = 0.0007085414921708109
$ws.Range("I7").Value = 0.3924995450689984
$ws.Range("J7").Value = 0.3924995450689983
$ws.Range("M7").Value = 15.47987166666667
$ws.Range("N7").Value = 46.439615
$ws.Range("O7").Value = 0.7960757698994193
$ws.Range("P7").Value = 0.7960757698994194
$ws.Range("Q7").Value = 3669.842422654203
$ws.Range("R7").Value = 33028.58180388783
$ws.Range("S7").Value = 0.3124593775259747
$ws.Range("T7").Value = 0.3124593775259747
$ws.Range("I8").Value = 0.3924995450689984
$ws.Range("J8").Value = 0.3924995450689983
$ws.Range("O8").Value = 0.05100527512565552
$ws.Range("P8").Value = 0.05100527512565553
$ws.Range("S8").Value = 0.02001954728293889
$ws.Range("T8").Value = 0.02001954728293889
$ws.Range("I9").Value = 0.3924995450689984
$ws.Range("J9").Value = 0.3924995450689983
$ws.Range("M9").Value = 0.171678
$ws.Range("N9").Value = 0.515034
$ws.Range("O9").Value = 0.008828800326496623
$ws.Range("P9").Value = 0.008828800326496624
$ws.Range("Q9").Value = 40.700027816968
$ws.Range("R9").Value = 366.300250352712
$ws.Range("S9").Value = 0.003465300111654949
$ws.Range("T9").Value = 0.003465300111654949
$ws.Range("I10").Value = 0.3924995450689984
$ws.Range("J10").Value = 0.3924995450689983
$ws.Range("M10").Value = 2.688466333333334
$ws.Range("N10").Value = 8.065399000000001
$ws.Range("O10").Value = 0.1382584398787761
$ws.Range("P10").Value = 0.1382584398787761
$ws.Range("Q10").Value = 637.3597930523927
$ws.Range("R10").Value = 5736.238137471533
$ws.Range("S10").Value = 0.05426637475436907
$ws.Range("T10").Value = 0.05426637475436906
$ws.Range("I11").Value = 0.3924995450689984
$ws.Range("J11").Value = 0.3924995450689983
$ws.Range("M11").Value = 0.113399
$ws.Range("N11").Value = 0.340197
$ws.Range("O11").Value = 0.005831714769652435
$ws.Range("P11").Value = 0.005831714769652436
$ws.Range("Q11").Value = 26.88371517851067
$ws.Range("R11").Value = 241.953436606596
$ws.Range("S11").Value = 0.002288945394060739
$ws.Range("T11").Value = 0.002288945394060739
$ws.Range("G12").Value = 138.1628113333333
$ws.Range("H12").Value = 414.488434
$ws.Range("I12").Value = 0.2287442882675098
$ws.Range("J12").Value = 0.2287442882675098
$ws.Range("M12").Value = 15.47987166666667
$ws.Range("N12").Value = 46.439615
$ws.Range("O12").Value = 0.7960757698994193
$ws.Range("P12").Value = 0.7960757698994194
$ws.Range("Q12").Value = 2138.742588545879
$ws.Range("R12").Value = 19248.68329691291
$ws.Range("S12").Value = 0.1820977853926526
$ws.Range("T12").Value = 0.1820977853926526
$ws.Range("G13").Value = 138.1628113333333
$ws.Range("H13").Value = 414.488434
$ws.Range("I13").Value = 0.2287442882675098
$ws.Range("J13").Value = 0.2287442882675098
$ws.Range("O13").Value = 0.05100527512565552
$ws.Range("P13").Value = 0.05100527512565553
$ws.Range("Q13").Value = 137.031119745702
$ws.Range("R13").Value = 1233.280077711318
$ws.Range("S13").Value = 0.0116671653565066
$ws.Range("T13").Value = 0.0116671653565066
$ws.Range("G14").Value = 138.1628113333333
$ws.Range("H14").Value = 414.488434
$ws.Range("I14").Value = 0.2287442882675098
$ws.Range("J14").Value = 0.2287442882675098
$ws.Range("M14").Value = 0.171678
$ws.Range("N14").Value = 0.515034
$ws.Range("O14").Value = 0.008828800326496623
$ws.Range("P14").Value = 0.008828800326496624
$ws.Range("Q14").Value = 23.719515124084
$ws.Range("R14").Value = 213.475636116756
$ws.Range("S14").Value = 0.002019537646940428
$ws.Range("T14").Value = 0.002019537646940428
$ws.Range("G15").Value = 138.1628113333333
$ws.Range("H15").Value = 414.488434
$ws.Range("I15").Value = 0.2287442882675098
$ws.Range("J15").Value = 0.2287442882675098
$ws.Range("M15").Value = 2.688466333333334
$ws.Range("N15").Value = 8.065399000000001
$ws.Range("O15").Value = 0.1382584398787761
$ws.Range("P15").Value = 0.1382584398787761
$ws.Range("Q15").Value = 371.4460667883519
$ws.Range("R15").Value = 3343.014601095166
$ws.Range("S15").Value = 0.03162582842704693
$ws.Range("T15").Value = 0.03162582842704692
$ws.Range("G16").Value = 138.1628113333333
$ws.Range("H16").Value = 414.488434
$ws.Range("I16").Value = 0.2287442882675098
$ws.Range("J16").Value = 0.2287442882675098
$ws.Range("M16").Value = 0.113399
$ws.Range("N16").Value = 0.340197
$ws.Range("O16").Value = 0.005831714769652435
$ws.Range("P16").Value = 0.005831714769652436
$ws.Range("Q16").Value = 15.66752464238867
$ws.Range("R16").Value = 141.007721781498
$ws.Range("S16").Value = 0.001333971444363271
$ws.Range("T16").Value = 0.001333971444363271
$ws.Range("G17").Value = 49.051656
$ws.Range("H17").Value = 147.154968
$ws.Range("I17").Value = 0.08121060965524597
$ws.Range("J17").Value = 0.08121060965524596
$ws.Range("M17").Value = 15.47987166666667
$ws.Range("N17").Value = 46.439615
$ws.Range("O17").Value = 0.7960757698994193
$ws.Range("P17").Value = 0.7960757698994194
$ws.Range("Q17").Value = 759.3133399174801
$ws.Range("R17").Value = 6833.82005925732
$ws.Range("S17").Value = 0.06464979860530115
$ws.Range("T17").Value = 0.06464979860530115
$ws.Range("G18").Value = 49.051656
$ws.Range("H18").Value = 147.154968
$ws.Range("I18").Value = 0.08121060965524597
$ws.Range("J18").Value = 0.08121060965524596
$ws.Range("O18").Value = 0.05100527512565552
$ws.Range("P18").Value = 0.05100527512565553
$ws.Range("Q18").Value = 48.649873885704
$ws.Range("R18").Value = 437.848864971336
$ws.Range("S18").Value = 0.004142169488588037
$ws.Range("T18").Value = 0.004142169488588037
$ws.Range("G19").Value = 49.051656
$ws.Range("H19").Value = 147.154968
$ws.Range("I19").Value = 0.08121060965524597
$ws.Range("J19").Value = 0.08121060965524596
$ws.Range("M19").Value = 0.171678
$ws.Range("N19").Value = 0.515034
$ws.Range("O19").Value = 0.008828800326496623
$ws.Range("P19").Value = 0.008828800326496624
$ws.Range("Q19").Value = 8.421090198768001
$ws.Range("R19").Value = 75.78981178891199
$ws.Range("S19").Value = 0.0007169922570392254
$ws.Range("T19").Value = 0.0007169922570392255
$ws.Range("G20").Value = 49.051656
$ws.Range("H20").Value = 147.154968
$ws.Range("I20").Value = 0.08121060965524597
$ws.Range("J20").Value = 0.08121060965524596
$ws.Range("M20").Value = 2.688466333333334
$ws.Range("N20").Value = 8.065399000000001
$ws.Range("O20").Value = 0.1382584398787761
$ws.Range("P20").Value = 0.1382584398787761
$ws.Range("Q20").Value = 131.873725750248
$ws.Range("R20").Value = 1186.863531752232
$ws.Range("S20").Value = 0.01122805219253858
$ws.Range("T20").Value = 0.01122805219253857
$ws.Range("G21").Value = 49.051656
$ws.Range("H21").Value = 147.154968
$ws.Range("I21").Value = 0.08121060965524597
$ws.Range("J21").Value = 0.08121060965524596
$ws.Range("M21").Value = 0.113399
$ws.Range("N21").Value = 0.340197
$ws.Range("O21").Value = 0.005831714769652435
$ws.Range("P21").Value = 0.005831714769652436
$ws.Range("Q21").Value = 5.562408738744001
$ws.Range("R21").Value = 50.061678648696
$ws.Range("S21").Value = 0.0004735971117789766
$ws.Range("T21").Value = 0.0004735971117789766
$ws.Range("G22").Value = 106.3337146666667
$ws.Range("H22").Value = 319.001144
$ws.Range("I22").Value = 0.1760475894022206
$ws.Range("J22").Value = 0.1760475894022206
$ws.Range("M22").Value = 15.47987166666667
$ws.Range("N22").Value = 46.439615
$ws.Range("O22").Value = 0.7960757698994193
$ws.Range("P22").Value = 0.7960757698994194
$ws.Range("Q22").Value = 1646.032256879951
$ws.Range("R22").Value = 14814.29031191956
$ws.Range("S22").Value = 0.1401472202723096
$ws.Range("T22").Value = 0.1401472202723096
$ws.Range("G23").Value = 106.3337146666667
$ws.Range("H23").Value = 319.001144
$ws.Range("I23").Value = 0.1760475894022206
$ws.Range("J23").Value = 0.1760475894022206
$ws.Range("O23").Value = 0.05100527512565552
$ws.Range("P23").Value = 0.05100527512565553
$ws.Range("Q23").Value = 105.462735209832
$ws.Range("R23").Value = 949.164616888488
$ws.Range("S23").Value = 0.0089793557326687
$ws.Range("T23").Value = 0.0089793557326687
$ws.Range("G24").Value = 106.3337146666667
$ws.Range("H24").Value = 319.001144
$ws.Range("I24").Value = 0.1760475894022206
$ws.Range("J24").Value = 0.1760475894022206
$ws.Range("M24").Value = 0.171678
$ws.Range("N24").Value = 0.515034
$ws.Range("O24").Value = 0.008828800326496623
$ws.Range("P24").Value = 0.008828800326496624
$ws.Range("Q24").Value = 18.255159466544
$ws.Range("R24").Value = 164.296435198896
$ws.Range("S24").Value = 0.001554289014793269
$ws.Range("T24").Value = 0.001554289014793269
$ws.Range("G25").Value = 106.3337146666667
$ws.Range("H25").Value = 319.001144
$ws.Range("I25").Value = 0.1760475894022206
$ws.Range("J25").Value = 0.1760475894022206
$ws.Range("M25").Value = 2.688466333333334
$ws.Range("N25").Value = 8.065399000000001
$ws.Range("O25").Value = 0.1382584398787761
$ws.Range("P25").Value = 0.1382584398787761
$ws.Range("Q25").Value = 285.8746119796063
$ws.Range("R25").Value = 2572.871507816456
$ws.Range("S25").Value = 0.02434006505517037
$ws.Range("T25").Value = 0.02434006505517037
$ws.Range("G26").Value = 106.3337146666667
$ws.Range("H26").Value = 319.001144
$ws.Range("I26").Value = 0.1760475894022206
$ws.Range("J26").Value = 0.1760475894022206
$ws.Range("M26").Value = 0.113399
$ws.Range("N26").Value = 0.340197
$ws.Range("O26").Value = 0.005831714769652435
$ws.Range("P26").Value = 0.005831714769652436
$ws.Range("Q26").Value = 12.05813690948533
$ws.Range("R26").Value = 108.523232185368
$ws.Range("S26").Value = 0.0007085414921708108
$ws.Range("T26").Value = 0.0007085414921708109
